$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.939.81"
$ws.Range("E2").Value = "  -3.63%  "
$ws.Range("D3").Value = "3.062.54"
$ws.Range("E3").Value = "  -2.90%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'542.91"
$ws.Range("E5").Value = "  -4.52%  "
$ws.Range("D6").Value = "'134.26"
$ws.Range("E6").Value = "  -9.76%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.051.84"
$ws.Range("E8").Value = "  -3.00%  "
$ws.Range("D9").Value = "'0.487"
$ws.Range("E9").Value = "  -2.53%  "
$ws.Range("D10").Value = "'6.47"
$ws.Range("E10").Value = "  -6.93%  "
$ws.Range("D11").Value = "'0.154"
$ws.Range("E11").Value = "  -2.56%  "
$ws.Range("D12").Value = "'0.454"
$ws.Range("E12").Value = "  -2.31%  "
$ws.Range("D13").Value = "'34.39"
$ws.Range("E13").Value = "  -4.61%  "
$ws.Range("D14").Value = "'0.0000214"
$ws.Range("E14").Value = "  -4.06%  "
$ws.Range("D15").Value = "3.556.68"
$ws.Range("E15").Value = "  -2.84%  "
$ws.Range("D16").Value = "63.037.80"
$ws.Range("E16").Value = "  -3.44%  "
$ws.Range("E17").Value = "  -1.68%  "
$ws.Range("D18").Value = "3.070.52"
$ws.Range("E18").Value = "  -2.60%  "
$ws.Range("D19").Value = "'6.56"
$ws.Range("E19").Value = "  -2.99%  "
$ws.Range("D20").Value = "'480.44"
$ws.Range("E20").Value = "  -8.73%  "
$ws.Range("D21").Value = "'13.26"
$ws.Range("E21").Value = "  -4.52%  "
$ws.Range("D22").Value = "'0.691"
$ws.Range("E22").Value = "  -1.99%  "
$ws.Range("D23").Value = "'7.05"
$ws.Range("E23").Value = "  -5.57%  "
$ws.Range("D24").Value = "'77.11"
$ws.Range("E24").Value = "  -2.22%  "
$ws.Range("D25").Value = "'12.01"
$ws.Range("E25").Value = "  -6.07%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'2.69"
$ws.Range("E27").Value = "  -4.53%  "
$ws.Range("D28").Value = "'8.12"
$ws.Range("E28").Value = "  -7.61%  "
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").Value = "'1.91"
$ws.Range("E30").Value = "  -10.49%  "
$ws.Range("D31").Value = "'26.04"
$ws.Range("E31").Value = "  -1.23%  "
$ws.Range("B32").Value = "OKB"
$ws.Range("C32").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D32").Value = "'60.93"
$ws.Range("E32").Value = "  +14.39%  "
$ws.Range("B33").Value = "Mantle"
$ws.Range("C33").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D33").Value = "'1.12"
$ws.Range("E33").Value = "  -2.24%  "
$ws.Range("D34").Value = "'2.46"
$ws.Range("E34").Value = "  -8.49%  "
$ws.Range("D35").Value = "'522.93"
$ws.Range("E35").Value = "  -5.94%  "
$ws.Range("D36").Value = "'5.85"
$ws.Range("E36").Value = "  -3.77%  "
$ws.Range("D37").Value = "'5.07"
$ws.Range("E37").Value = "  -7.16%  "
$ws.Range("D38").Value = "'0.0395"
$ws.Range("E38").Value = "  -10.80%  "
$ws.Range("D39").Value = "3.061.33"
$ws.Range("E39").Value = "  -0.74%  "
$ws.Range("D40").Value = "'0.0779"
$ws.Range("E40").Value = "  -5.43%  "
$ws.Range("D41").Value = "'0.117"
$ws.Range("E41").Value = "  -4.01%  "
$ws.Range("D42").Value = "'7.98"
$ws.Range("E42").Value = "  -3.47%  "
$ws.Range("D43").Value = "'2.61"
$ws.Range("E43").Value = "  -9.09%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.249"
$ws.Range("E45").Value = "  -3.35%  "
$ws.Range("D46").Value = "'2.01"
$ws.Range("E46").Value = "  -9.09%  "
$ws.Range("D47").Value = "'120.47"
$ws.Range("E47").Value = "  +1.51%  "
$ws.Range("D48").Value = "'23.92"
$ws.Range("E48").Value = "  -4.64%  "
$ws.Range("D49").Value = "'0.106"
$ws.Range("E49").Value = "  -2.95%  "
$ws.Range("D50").Value = "0.0₃0495"
$ws.Range("E50").Value = "  -5.87%  "
$ws.Range("D51").Value = "'2.33"
$ws.Range("E51").Value = "  +58.72%  "
